# Update "想去人数" (want-to-go count, column F) and "最低票价" (lowest
# ticket price, column G) figures on the "展览" and "全部类型" sheets to
# reflect refreshed scrape data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("G2").Value = 55
$ws.Range("G4").Value = 55
$ws.Range("F5").Value = 6764
$ws.Range("F6").Value = 87
$ws.Range("F7").Value = 9
$ws.Range("F8").Value = 437
$ws.Range("F10").Value = 6343
$ws.Range("G13").Value = 68
$ws.Range("F17").Value = 126
$ws.Range("F22").Value = 4700
$ws.Range("F23").Value = 70
$ws.Range("F25").Value = 154
$ws.Range("F27").Value = 99

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G2").Value = 55
$ws.Range("G4").Value = 55
$ws.Range("F5").Value = 6764
$ws.Range("F6").Value = 87
$ws.Range("F7").Value = 9
$ws.Range("F8").Value = 437
$ws.Range("F10").Value = 6343
$ws.Range("G13").Value = 68
$ws.Range("F17").Value = 126
$ws.Range("F22").Value = 4700
$ws.Range("F24").Value = 70
$ws.Range("F26").Value = 154
$ws.Range("F28").Value = 99
